# Updates the cryptocurrency price table (columns D/E) with refreshed market
# data, and swaps the EnergySwap/Aptos rows (47/48) per the upstream source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cell reference, new text value) pairs taken from the refreshed feed.
$updates = @(
    @("D2", "29.986.24"),
    @("E2", "  -0.55%  "),
    @("D3", "1.897.39"),
    @("E3", "  -0.98%  "),
    @("D4", "1.002"),
    @("E4", "  +0.13%  "),
    @("D5", "0.8280"),
    @("E5", "  +4.17%  "),
    @("D6", "241.78"),
    @("E6", "  -0.52%  "),
    @("D7", "1.002"),
    @("E7", "  +0.20%  "),
    @("D8", "0.3271"),
    @("E8", "  +2.59%  "),
    @("D9", "26.46"),
    @("E9", "  +0.10%  "),
    @("D10", "0.07028"),
    @("E10", "  +0.80%  "),
    @("D11", "0.08085"),
    @("E11", "  +0.76%  "),
    @("D12", "0.7590"),
    @("E12", "  +0.68%  "),
    @("D13", "1.905.83"),
    @("E13", "  -0.52%  "),
    @("D14", "5.234"),
    @("E14", "  -0.02%  "),
    @("D15", "92.09"),
    @("E15", "  -1.73%  "),
    @("D16", "29.984.17"),
    @("E16", "  -0.61%  "),
    @("D17", "14.08"),
    @("E17", "  +0.08%  "),
    @("D18", "5.846"),
    @("E18", "  -2.80%  "),
    @("D19", "243.90"),
    @("E19", "  -2.31%  "),
    @("D20", "0.000007737"),
    @("E20", "  -1.27%  "),
    @("D21", "1.002"),
    @("E21", "  +0.17%  "),
    @("D22", "2.147.50"),
    @("E22", "  -0.86%  "),
    @("D23", "1.002"),
    @("E23", "  +0.17%  "),
    @("D24", "6.946"),
    @("E24", "  -0.70%  "),
    @("D25", "0.1734"),
    @("E25", "  +23.84%  "),
    @("D26", "9.246"),
    @("E26", "  -1.12%  "),
    @("D27", "165.61"),
    @("E27", "  -2.14%  "),
    @("D28", "18.88"),
    @("E28", "  -0.76%  "),
    @("D29", "2.087"),
    @("E29", "  +1.39%  "),
    @("D30", "1.362"),
    @("D31", "1.515"),
    @("E31", "  -0.92%  "),
    @("D32", "0.05914"),
    @("E32", "  +9.65%  "),
    @("D33", "4.270"),
    @("E33", "  -2.40%  "),
    @("D34", "4.062"),
    @("E34", "  -1.61%  "),
    @("D35", "1.262"),
    @("E35", "  -0.79%  "),
    @("D36", "0.7301"),
    @("E36", "  -1.43%  "),
    @("D37", "2.725"),
    @("E37", "  -0.19%  "),
    @("D38", "0.01914"),
    @("E38", "  -1.00%  "),
    @("E39", "  -0.67%  "),
    @("D40", "0.4431"),
    @("E40", "  -0.82%  "),
    @("D41", "72.24"),
    @("E41", "  -0.77%  "),
    @("D42", "5.849"),
    @("E42", "  -5.73%  "),
    @("D43", "0.8499"),
    @("E43", "  +1.76%  "),
    @("D44", "1.002"),
    @("E44", "  +0.15%  "),
    @("D45", "1.892"),
    @("E45", "  -0.83%  "),
    @("D46", "101.91"),
    @("E46", "  +1.19%  "),
    @("B47", "Aptos"),
    @("C47", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"),
    @("D47", "7.534"),
    @("E47", "  -1.29%  "),
    @("B48", "EnergySwap"),
    @("C48", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D48", "9.775"),
    @("E48", "  -1.10%  "),
    @("D49", "989.89"),
    @("E49", "  +2.32%  "),
    @("D50", "2.045.77"),
    @("E50", "  -0.79%  "),
    @("D51", "1.512"),
    @("E51", "  -0.24%  "),
)

foreach ($pair in $updates) {
    $ref = $pair[0]
    $val = $pair[1]
    $cell = $ws.Range($ref)
    if ($ref -match "^D" -and $val -match "^[0-9.]+$") {
        # Values in column D are textual price strings (some use "."
        # as a thousands separator) that would otherwise be auto-coerced
        # into numbers by plain Value assignment. Force text, write, then
        # restore the default "Normal" style so no stray formatting sticks.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
